# disabled_students_report.xlsx: two newly-disabled students got logged, and
# the first row's disabling admin/reason/photo/timestamp were corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: fix up the existing record (disabled_by / reason / image_path / date_disabled) ---
$ws.Range("E2").Value = "Dennis"
$ws.Range("F2").Value = "ddddd"
# no photo was attached this time -> blank (but still text-typed) cell, same
# as every other empty image_path cell this export produces
$ws.Range("G2").Value = "'"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = "2025-05-25 15:07:32"

# --- Seed rows 3 & 4 with the repeated student info by copying row 2's
#     B:D cells (keeps them text-typed like the source data, e.g. "28527"
#     stays a string instead of being re-inferred as a number) ---
$ws.Range("B2:D2").Copy()
$ws.Range("B3:D3").PasteSpecial()
$ws.Range("B2:D2").Copy()
$ws.Range("B4:D4").PasteSpecial()

# --- Row 3: new disabled-student record ---
$ws.Range("A3").Value = 2
$ws.Range("E3").Value = "Dennis"
$ws.Range("F3").Value = "dddd"
$ws.Range("G3").Value = "SBoys print25043012510.jpg"
$ws.Range("H3").Value = "2025-05-25 15:07:59"

# --- Row 4: new disabled-student record ---
$ws.Range("A4").Value = 3
$ws.Range("E4").Value = "Mellisa"
$ws.Range("F4").Value = "ewwewe"
$ws.Range("G4").Value = "'"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = "2025-05-25 15:13:53"
